$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for account 004862672 / RENATO / 104509.27 (originally row 2)
$ws.Rows.Item(2).Delete()

# Replace the row for account 005681354 / MATHEUS / 8000 (now row 4, after the
# deletion above shifted everything up by one) with the new record for
# 004983378 / MARCELO / 1767.28.
# Force text format temporarily so the leading zeros in the account number
# survive (Excel would otherwise auto-convert the numeric-looking string to
# a number), then clear the formatting again so the cell's style matches the
# plain (unstyled) look of every other data cell in the column.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "004983378"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "MARCELO"
$ws.Range("C4").Value = 1767.28
